# Updated figure with kinesis
#
# 1) The cached "datetimeFigureOut" footer-date text (1/23/14 -> 9/6/14) is
#    stored on every slide layout's and the slide master's "Date Placeholder"
#    shape. Update each one explicitly (by its known shape index) so we don't
#    need COM-call-heavy loops.
# 2) On slide 1, the streaming-sources diagram gets two label edits inside a
#    deeply nested group: "HDFS" -> "HDFS/S3" and "ZeroMQ" -> "Kinesis".

$p = $ppt.ActivePresentation

# --- 1. Refresh the cached date text everywhere it is cached -------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$layouts = $master.CustomLayouts

$master.Shapes.Item(3).TextFrame.TextRange.Text = "9/6/14"

$layouts.Item(1).Shapes.Item(3).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(2).Shapes.Item(3).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(3).Shapes.Item(3).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(4).Shapes.Item(4).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(5).Shapes.Item(6).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(6).Shapes.Item(2).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(8).Shapes.Item(4).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(9).Shapes.Item(4).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(10).Shapes.Item(3).TextFrame.TextRange.Text = "9/6/14"
$layouts.Item(11).Shapes.Item(3).TextFrame.TextRange.Text = "9/6/14"

# --- 2. Update the two storage/ingestion labels on slide 1 ---------------
$slide1 = $p.Slides.Item(1)
$diagram = $slide1.Shapes.Item(1)

# The diagram is one big nested group; GroupItems flattens it in document
# order, so item 8 is the "HDFS" rounded rectangle and item 9 is the
# "ZeroMQ" rounded rectangle right after it.
$hdfsShape = $diagram.GroupItems.Item(8)
$queueShape = $diagram.GroupItems.Item(9)

$hdfsShape.TextFrame.TextRange.Text = "HDFS/S3"
$queueShape.TextFrame.TextRange.Text = "Kinesis"
